$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.006840989660467756
$ws.Range("C2").Value = -0.0006251249918216573

$ws.Range("B3").Value = 0.01639476845139143
$ws.Range("C3").Value = 0.004728100486627085

$ws.Range("C4").Value = 0.127310736114822

$ws.Range("B5").Value = -0.006635612259827894
$ws.Range("C5").Value = -0.0006063577443517487

$ws.Range("B6").Value = -0.2154074384393425
$ws.Range("C6").Value = -0.01968378210892752

$ws.Range("B7").Value = -0.08974641073382372
$ws.Range("C7").Value = -0.002853739048248372
